# APPL should always have been AAPL
# Rename the placeholder ticker symbols "A" and "B" (used in the
# buy_orders / sell_orders test fixtures) to "ALOTOFLOVE" and
# "BERIGHTBACK" respectively, and leave the sheet selection on the
# "symbol" cell of the first lot in each of those two sheets.

$wb = $excel.ActiveWorkbook

$buyOrders = $wb.Worksheets.Item("buy_orders")
$buyOrders.Range("B2").Value = "ALOTOFLOVE"
$buyOrders.Range("B3").Value = "BERIGHTBACK"

$sellOrders = $wb.Worksheets.Item("sell_orders")
$sellOrders.Range("B2").Value = "ALOTOFLOVE"

# Select B2 on sell_orders first ...
[void]$sellOrders.Select()
[void]$sellOrders.Range("B2").Select()

# ... then leave buy_orders as the active / front-most sheet with B2
# selected, matching the saved workbook's view state.
[void]$buyOrders.Select()
[void]$buyOrders.Range("B2").Select()
